# ============================================================
# Update master qubit with today's (2019-10-28) run data.
# Adds 24 new sample rows (447-470) to the qubit-iso sheet,
# mirroring the structure/formulas of the previous batch
# (rows 441-446), and updates the sheet view selection.
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 447
$lastRow  = 470
$nRows    = $lastRow - $firstRow + 1

# ------------------------------------------------------------
# 1) Apply cell formatting for the new rows by copying it from
#    existing rows that already carry the desired per-column
#    styles, so we reuse the workbook's existing style entries
#    instead of creating new ones.
# ------------------------------------------------------------

# D (date), I, J, K, L, M, N reuse the formatting already used
# on the prior batch's rows (row 441). Apply this first since it
# spans the full D:N block.
$ws.Range("D441:N441").Copy() | Out-Null
$ws.Range("D$firstRow`:N$lastRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# A:C and E:H use the explicit-black-font style also seen on
# e.g. row 313 (style index "1" in the original file). Apply this
# afterwards so it isn't clobbered by the D:N paste above.
$ws.Range("A313:C313").Copy() | Out-Null
$ws.Range("A$firstRow`:C$lastRow").PasteSpecial(-4122) | Out-Null

$ws.Range("E313:H313").Copy() | Out-Null
$ws.Range("E$firstRow`:H$lastRow").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ------------------------------------------------------------
# 2) Write the literal data values for the 24 new rows.
# ------------------------------------------------------------

# Columns A,B,C,D,E,F,G,H,I
$data1 = New-Object 'object[,]' 24,9
$data1[0,0] = "2019-10-28_171734"; $data1[0,1] = "RNA High sensitivity"; $data1[0,2] = "Sample_#191028-172129"; $data1[0,3] = 43766.723252314812; $data1[0,4] = 244; $data1[0,5] = 24.4; $data1[0,6] = 2; $data1[0,7] = 100; $data1[0,8] = 282
$data1[1,0] = "2019-10-28_171734"; $data1[1,1] = "RNA High sensitivity"; $data1[1,2] = "Sample_#191028-172120"; $data1[1,3] = 43766.723252314812; $data1[1,4] = 152; $data1[1,5] = 15.2; $data1[1,6] = 2; $data1[1,7] = 100; $data1[1,8] = 268
$data1[2,0] = "2019-10-28_171734"; $data1[2,1] = "RNA High sensitivity"; $data1[2,2] = "Sample_#191028-172111"; $data1[2,3] = 43766.723252314812; $data1[2,4] = 183; $data1[2,5] = 18.3; $data1[2,6] = 2; $data1[2,7] = 100; $data1[2,8] = 374
$data1[3,0] = "2019-10-28_171734"; $data1[3,1] = "RNA High sensitivity"; $data1[3,2] = "Sample_#191028-172102"; $data1[3,3] = 43766.723252314812; $data1[3,4] = 231; $data1[3,5] = 23.1; $data1[3,6] = 2; $data1[3,7] = 100; $data1[3,8] = 281
$data1[4,0] = "2019-10-28_171734"; $data1[4,1] = "RNA High sensitivity"; $data1[4,2] = "Sample_#191028-172054"; $data1[4,3] = 43766.723252314812; $data1[4,4] = 84.9; $data1[4,5] = 8.49; $data1[4,6] = 2; $data1[4,7] = 100; $data1[4,8] = 296
$data1[5,0] = "2019-10-28_171734"; $data1[5,1] = "RNA High sensitivity"; $data1[5,2] = "Sample_#191028-172044"; $data1[5,3] = 43766.723252314812; $data1[5,4] = 84.7; $data1[5,5] = 8.47; $data1[5,6] = 2; $data1[5,7] = 100; $data1[5,8] = 376
$data1[6,0] = "2019-10-28_171734"; $data1[6,1] = "RNA High sensitivity"; $data1[6,2] = "Sample_#191028-172035"; $data1[6,3] = 43766.723252314812; $data1[6,4] = 52.3; $data1[6,5] = 5.23; $data1[6,6] = 2; $data1[6,7] = 100; $data1[6,8] = 294
$data1[7,0] = "2019-10-28_171734"; $data1[7,1] = "RNA High sensitivity"; $data1[7,2] = "Sample_#191028-172027"; $data1[7,3] = 43766.723252314812; $data1[7,4] = 355; $data1[7,5] = 35.5; $data1[7,6] = 2; $data1[7,7] = 100; $data1[7,8] = 279
$data1[8,0] = "2019-10-28_171734"; $data1[8,1] = "RNA High sensitivity"; $data1[8,2] = "Sample_#191028-172018"; $data1[8,3] = 43766.723252314812; $data1[8,4] = 223; $data1[8,5] = 22.3; $data1[8,6] = 2; $data1[8,7] = 100; $data1[8,8] = 377
$data1[9,0] = "2019-10-28_171734"; $data1[9,1] = "RNA High sensitivity"; $data1[9,2] = "Sample_#191028-172009"; $data1[9,3] = 43766.723252314812; $data1[9,4] = 23; $data1[9,5] = 2.3; $data1[9,6] = 2; $data1[9,7] = 100; $data1[9,8] = 243
$data1[10,0] = "2019-10-28_171734"; $data1[10,1] = "RNA High sensitivity"; $data1[10,2] = "Sample_#191028-172000"; $data1[10,3] = 43766.723252314812; $data1[10,4] = 262; $data1[10,5] = 26.2; $data1[10,6] = 2; $data1[10,7] = 100; $data1[10,8] = 226
$data1[11,0] = "2019-10-28_171734"; $data1[11,1] = "RNA High sensitivity"; $data1[11,2] = "Sample_#191028-171952"; $data1[11,3] = 43766.723252314812; $data1[11,4] = 340; $data1[11,5] = 34; $data1[11,6] = 2; $data1[11,7] = 100; $data1[11,8] = 213
$data1[12,0] = "2019-10-28_171734"; $data1[12,1] = "RNA High sensitivity"; $data1[12,2] = "Sample_#191028-171943"; $data1[12,3] = 43766.723252314812; $data1[12,4] = 180; $data1[12,5] = 18; $data1[12,6] = 2; $data1[12,7] = 100; $data1[12,8] = 227
$data1[13,0] = "2019-10-28_171734"; $data1[13,1] = "RNA High sensitivity"; $data1[13,2] = "Sample_#191028-171933"; $data1[13,3] = 43766.723252314812; $data1[13,4] = 62.7; $data1[13,5] = 6.27; $data1[13,6] = 2; $data1[13,7] = 100; $data1[13,8] = 201
$data1[14,0] = "2019-10-28_171734"; $data1[14,1] = "RNA High sensitivity"; $data1[14,2] = "Sample_#191028-171924"; $data1[14,3] = 43766.723252314812; $data1[14,4] = 134; $data1[14,5] = 13.4; $data1[14,6] = 2; $data1[14,7] = 100; $data1[14,8] = 248
$data1[15,0] = "2019-10-28_171734"; $data1[15,1] = "RNA High sensitivity"; $data1[15,2] = "Sample_#191028-171916"; $data1[15,3] = 43766.723252314812; $data1[15,4] = 216; $data1[15,5] = 21.6; $data1[15,6] = 2; $data1[15,7] = 100; $data1[15,8] = 240
$data1[16,0] = "2019-10-28_171734"; $data1[16,1] = "RNA High sensitivity"; $data1[16,2] = "Sample_#191028-171907"; $data1[16,3] = 43766.723252314812; $data1[16,4] = 130; $data1[16,5] = 13; $data1[16,6] = 2; $data1[16,7] = 100; $data1[16,8] = 241
$data1[17,0] = "2019-10-28_171734"; $data1[17,1] = "RNA High sensitivity"; $data1[17,2] = "Sample_#191028-171859"; $data1[17,3] = 43766.723252314812; $data1[17,4] = 317; $data1[17,5] = 31.7; $data1[17,6] = 2; $data1[17,7] = 100; $data1[17,8] = 259
$data1[18,0] = "2019-10-28_171734"; $data1[18,1] = "RNA High sensitivity"; $data1[18,2] = "Sample_#191028-171850"; $data1[18,3] = 43766.723252314812; $data1[18,4] = 273; $data1[18,5] = 27.3; $data1[18,6] = 2; $data1[18,7] = 100; $data1[18,8] = 310
$data1[19,0] = "2019-10-28_171734"; $data1[19,1] = "RNA High sensitivity"; $data1[19,2] = "Sample_#191028-171842"; $data1[19,3] = 43766.723252314812; $data1[19,4] = 178; $data1[19,5] = 17.8; $data1[19,6] = 2; $data1[19,7] = 100; $data1[19,8] = 315
$data1[20,0] = "2019-10-28_171734"; $data1[20,1] = "RNA High sensitivity"; $data1[20,2] = "Sample_#191028-171833"; $data1[20,3] = 43766.723252314812; $data1[20,4] = 29; $data1[20,5] = 2.9; $data1[20,6] = 2; $data1[20,7] = 100; $data1[20,8] = 329
$data1[21,0] = "2019-10-28_171734"; $data1[21,1] = "RNA High sensitivity"; $data1[21,2] = "Sample_#191028-171824"; $data1[21,3] = 43766.723252314812; $data1[21,4] = 274; $data1[21,5] = 27.4; $data1[21,6] = 2; $data1[21,7] = 100; $data1[21,8] = 303
$data1[22,0] = "2019-10-28_171734"; $data1[22,1] = "RNA High sensitivity"; $data1[22,2] = "Sample_#191028-171816"; $data1[22,3] = 43766.723252314812; $data1[22,4] = 311; $data1[22,5] = 31.1; $data1[22,6] = 2; $data1[22,7] = 100; $data1[22,8] = 325
$data1[23,0] = "2019-10-28_171734"; $data1[23,1] = "RNA High sensitivity"; $data1[23,2] = "Sample_#191028-171806"; $data1[23,3] = 43766.723252314812; $data1[23,4] = 289; $data1[23,5] = 28.9; $data1[23,6] = 2; $data1[23,7] = 100; $data1[23,8] = 301
$ws.Range("A$firstRow`:I$lastRow").Value = $data1

# Columns K,L
$data2 = New-Object 'object[,]' 24,2
$data2[0,0] = 35; $data2[0,1] = 15
$data2[1,0] = 35; $data2[1,1] = 15
$data2[2,0] = 35; $data2[2,1] = 15
$data2[3,0] = 35; $data2[3,1] = 15
$data2[4,0] = 35; $data2[4,1] = 15
$data2[5,0] = 35; $data2[5,1] = 15
$data2[6,0] = 35; $data2[6,1] = 15
$data2[7,0] = 35; $data2[7,1] = 15
$data2[8,0] = 35; $data2[8,1] = 15
$data2[9,0] = 35; $data2[9,1] = 15
$data2[10,0] = 35; $data2[10,1] = 15
$data2[11,0] = 35; $data2[11,1] = 15
$data2[12,0] = 35; $data2[12,1] = 15
$data2[13,0] = 35; $data2[13,1] = 15
$data2[14,0] = 35; $data2[14,1] = 15
$data2[15,0] = 35; $data2[15,1] = 15
$data2[16,0] = 35; $data2[16,1] = 15
$data2[17,0] = 35; $data2[17,1] = 15
$data2[18,0] = 35; $data2[18,1] = 15
$data2[19,0] = 35; $data2[19,1] = 15
$data2[20,0] = 35; $data2[20,1] = 15
$data2[21,0] = 35; $data2[21,1] = 15
$data2[22,0] = 35; $data2[22,1] = 15
$data2[23,0] = 35; $data2[23,1] = 15
$ws.Range("K$firstRow`:L$lastRow").Value = $data2

# Column J (extraction_method) and N (sample_type) are constant for the batch
$ws.Range("J$firstRow`:J$lastRow").Value = "Zymo_microprep"
$ws.Range("N$firstRow`:N$lastRow").Value = "pellet"

# ------------------------------------------------------------
# 3) Column M: extend the shared '(F)*(L-G)' yield formula
# ------------------------------------------------------------
$ws.Range("M$firstRow`:M$lastRow").Formula = "=(F$firstRow)*(L$firstRow-G$firstRow)"

# ------------------------------------------------------------
# 4) Update the sheet view to reflect where the user ended up
# ------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 430
$win.ScrollColumn = 1
$ws.Range("P459").Select() | Out-Null
